$d = $word.ActiveDocument

$target = "　　　　　　　　　　院长：签名（盖章）"

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*$target*") {
        $p.Range.Delete()
        break
    }
}
